# lesson 71 - after homework - part 1
#
# Splits four runs of placeholder dots so that the (previously hidden)
# answer word is inserted in the middle of each blank, e.g.
#   "We have suffered a major …………………………."
# becomes
#   "We have suffered a major …" + "setback" + "………………………."
#
# The trick used to force Word to create a *new* <w:r> at an insertion
# point (rather than silently merging the inserted text into the
# neighbouring run, which happens whenever the formatting stays
# identical) is to toggle a character formatting property on the
# freshly inserted range immediately after inserting it: turning Bold
# on and back off again leaves the visible formatting unchanged but
# makes the engine split the run.

$d = $word.ActiveDocument

function Split-Insert($range, $splitPos, $text) {
    # $splitPos is an absolute document character offset; the insertion
    # happens right there, pushing the rest of the original run into a
    # new trailing run.
    $insertRange = $d.Range($splitPos, $splitPos)
    $insertRange.InsertAfter($text)
    $insertRange.Font.Bold = 1
    $insertRange.Font.Bold = 0
}

# ---------------------------------------------------------------------
# 1) "We have suffered a major ………………………….":
#    "We have suffered a major …" + "setback" + "………………………."
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("We have suffered a major ………………………….", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $r.Start
$splitPos = $start + ("We have suffered a major …").Length
Split-Insert $r $splitPos "setback"

# ---------------------------------------------------------------------
# 2) "Your achievement is remarkable , without any doubts ………………….":
#    "Your achievement is remarkable , without any doubts …" + "praiseworthy" + "………………."
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Your achievement is remarkable , without any doubts ………………….", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $r.Start
$splitPos = $start + ("Your achievement is remarkable , without any doubts …").Length
Split-Insert $r $splitPos "praiseworthy"

# ---------------------------------------------------------------------
# 3) "Today, the ………………………………in college job" (highlighted "shl" run):
#    "…" + "buzzword" + "……………………………"
#
#    This paragraph continues, after the dots, with several more runs
#    that all happen to share identical formatting ("in college job" +
#    " " + "hunting" + " " + 'is "information interview."'). Touching
#    the paragraph causes the engine to coalesce adjacent same-format
#    runs, so after inserting "buzzword" we must re-impose the original
#    run boundaries for the remainder of the paragraph (they are
#    untouched by the diff, so they must come out exactly as they were).
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Today, the ………………………………in college job", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $r.Start
$splitPos = $start + ("Today, the …").Length
Split-Insert $r $splitPos "buzzword"

$r2 = $d.Content
$r2.Find.Execute("in college job", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$icjStart = $r2.Start
$paraEnd = $r2.Paragraphs(1).Range.End

# original run boundaries, relative to the start of "in college job":
#   "in college job" | " " | "hunting" | " " | 'is "information interview."'
foreach ($boundary in @(14, 15, 22, 23)) {
    $tail = $d.Range($icjStart + $boundary, $paraEnd)
    $tail.Font.Bold = 1
    $tail.Font.Bold = 0
}

# ---------------------------------------------------------------------
# 4) "We ought to ………………………with that risky venture or else we may face flactuations":
#    "We ought to ………………………with that risky venture or else we m" + "ay face flu" + "ctuations"
#    (also fixes the typo "flactuations" -> "fluctuations")
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("We ought to ………………………with that risky venture or else we may face flactuations", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $r.Start
$p1 = $start + ("We ought to ………………………with that risky venture or else we m").Length
$p2 = $p1 + ("ay face fla").Length

$midRange = $d.Range($p1, $p2)
$midRange.Text = "ay face flu"
$midRange.Font.Bold = 1
$midRange.Font.Bold = 0
